$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2159311196958242
$ws.Range("B3").Value = 0.1913949238933697
$ws.Range("B4").Value = 0.1762869034279646
$ws.Range("B5").Value = 0.1701201473106835
$ws.Range("B6").Value = 0.1690955695415823
$ws.Range("B7").Value = 0.1762037765295901
$ws.Range("B8").Value = 0.2074802658001431
$ws.Range("B9").Value = 0.2684501888315936
$ws.Range("B10").Value = 0.3129948669844396
$ws.Range("B11").Value = 0.33319936175468
$ws.Range("B12").Value = 0.3408412057247006
$ws.Range("B13").Value = 0.3391958157679937
$ws.Range("B14").Value = 0.3338282482071406
$ws.Range("B15").Value = 0.3305392456628908
$ws.Range("B16").Value = 0.3116732098856971
$ws.Range("B17").Value = 0.3000838993125114
$ws.Range("B18").Value = 0.2934125106850161
$ws.Range("B19").Value = 0.2911527646049308
$ws.Range("B20").Value = 0.3013181776851468
$ws.Range("B21").Value = 0.3354050867553156
$ws.Range("B22").Value = 0.3576291578720543
$ws.Range("B23").Value = 0.3457728835225566
$ws.Range("B24").Value = 0.3007601869816199
$ws.Range("B25").Value = 0.2519980490937996
$ws.Range("D2").Value = 0.00328276805978156
$ws.Range("D3").Value = 0.003117985810526136
$ws.Range("D4").Value = 0.003018255499195988
$ws.Range("D5").Value = 0.002977983059874845
$ws.Range("D6").Value = 0.002971318247256605
$ws.Range("D7").Value = 0.003017710873152168
$ws.Range("D8").Value = 0.00322565360309035
$ws.Range("D9").Value = 0.003644731220816766
$ws.Range("D10").Value = 0.003959326331720803
$ws.Range("D11").Value = 0.004103863068408486
$ws.Range("D12").Value = 0.004158796915575635
$ws.Range("D13").Value = 0.004146957056267553
$ws.Range("D14").Value = 0.004108378505360832
$ws.Range("D15").Value = 0.00408477407517438
$ws.Range("D16").Value = 0.003949908852845851
$ws.Range("D17").Value = 0.003867535620365459
$ws.Range("D18").Value = 0.003820291157840217
$ws.Range("D19").Value = 0.003804318212985436
$ws.Range("D20").Value = 0.003876290507584201
$ws.Range("D21").Value = 0.004119704540567426
$ws.Range("D22").Value = 0.00427995902020939
$ws.Range("D23").Value = 0.00419432250900087
$ws.Range("D24").Value = 0.003872332071402695
$ws.Range("D25").Value = 0.003530170004244937
$ws.Range("E2").Value = 0.4354855266417559
$ws.Range("E3").Value = 0.3796236440278449
$ws.Range("E4").Value = 0.3454298728940586
$ws.Range("E5").Value = 0.3315199486717546
$ws.Range("E6").Value = 0.3292116194363217
$ws.Range("E7").Value = 0.3452421832880503
$ws.Range("E8").Value = 0.4162010468717057
$ws.Range("E9").Value = 0.5562976445450403
$ws.Range("E10").Value = 0.6599630079396519
$ws.Range("E11").Value = 0.7073182102084701
$ws.Range("E12").Value = 0.7252814645633805
$ws.Range("E13").Value = 0.7214113434755376
$ws.Range("E14").Value = 0.7087954249356869
$ws.Range("E15").Value = 0.7010719049055893
$ws.Range("E16").Value = 0.6568724169760856
$ws.Range("E17").Value = 0.629809763319642
$ws.Range("E18").Value = 0.6142624613427046
$ws.Range("E19").Value = 0.6090015003013463
$ws.Range("E20").Value = 0.6326887045236589
$ws.Range("E21").Value = 0.7125001682736638
$ws.Range("E22").Value = 0.7648426698322197
$ws.Range("E23").Value = 0.7368890890760866
$ws.Range("E24").Value = 0.6313870999754982
$ws.Range("E25").Value = 0.5182800985496812
$ws.Range("F2").Value = 0.4268455756567278
$ws.Range("F3").Value = 0.3933003553665202
$ws.Range("F4").Value = 0.3729420886923549
$ws.Range("F5").Value = 0.3647054950681934
$ws.Range("F6").Value = 0.3633413997772692
$ws.Range("F7").Value = 0.3728307665163442
$ws.Range("F8").Value = 0.4152293778963099
$ws.Range("F9").Value = 0.5002909339772401
$ws.Range("F10").Value = 0.5639956345128354
$ws.Range("F11").Value = 0.5932485116357213
$ws.Range("F12").Value = 0.604365721277091
$ws.Range("F13").Value = 0.6019696543234119
$ws.Range("F14").Value = 0.5941623311389179
$ws.Range("F15").Value = 0.5893853151532795
$ws.Range("F16").Value = 0.5620894304974513
$ws.Range("F17").Value = 0.5454146261146775
$ws.Range("F18").Value = 0.5358494138644119
$ws.Range("F19").Value = 0.5326151957723084
$ws.Range("F20").Value = 0.5471870247896504
$ws.Range("F21").Value = 0.5964544477321994
$ws.Range("F22").Value = 0.6288856518111032
$ws.Range("F23").Value = 0.6115551107154431
$ws.Range("F24").Value = 0.5463856568671446
$ws.Range("F25").Value = 0.4770700030462791
$ws.Range("G2").Value = 0.3367752227119354
$ws.Range("G3").Value = 0.3012008603001277
$ws.Range("G4").Value = 0.2794993337179079
$ws.Range("G5").Value = 0.2706907303260948
$ws.Range("G6").Value = 0.269230162170615
$ws.Range("G7").Value = 0.2793803970864275
$ws.Range("G8").Value = 0.3244794999739042
$ws.Range("G9").Value = 0.4140703990656505
$ws.Range("G10").Value = 0.4806440831150098
$ws.Range("G11").Value = 0.5111043360739984
$ws.Range("G12").Value = 0.5226648403903766
$ws.Range("G13").Value = 0.5201739184791165
$ws.Range("G14").Value = 0.5120549037030742
$ws.Range("G15").Value = 0.5070851600557091
$ws.Range("G16").Value = 0.4786570196494324
$ws.Range("G17").Value = 0.4612626808594769
$ws.Range("G18").Value = 0.4512744032264209
$ws.Range("G19").Value = 0.4478953566152768
$ws.Range("G20").Value = 0.4631126271809762
$ws.Range("G21").Value = 0.5144389481118878
$ws.Range("G22").Value = 0.5481347357242612
$ws.Range("G23").Value = 0.5301366292087835
$ws.Range("G24").Value = 0.4622762292460152
$ws.Range("G25").Value = 0.3897053291205168
$ws.Range("H2").Value = 0.3434920082859492
$ws.Range("H3").Value = 0.3306581927846537
$ws.Range("H4").Value = 0.3229880749787242
$ws.Range("H5").Value = 0.3199150623406126
$ws.Range("H6").Value = 0.3194079646665671
$ws.Range("H7").Value = 0.3229464183999653
$ws.Range("H8").Value = 0.3390232442816625
$ws.Range("H9").Value = 0.3722245405283218
$ws.Range("H10").Value = 0.3976543754028228
$ws.Range("H11").Value = 0.4094518639923876
$ws.Range("H12").Value = 0.4139524800987431
$ws.Range("H13").Value = 0.4129817139520924
$ws.Range("H14").Value = 0.409821466670877
$ws.Range("H15").Value = 0.4078900483665961
$ws.Range("H16").Value = 0.3968880122874339
$ws.Range("H17").Value = 0.3901974874650023
$ws.Range("H18").Value = 0.3863708487733675
$ws.Range("H19").Value = 0.3850789144428575
$ws.Range("H20").Value = 0.3909074710753089
$ws.Range("H21").Value = 0.4107488063632729
$ws.Range("H22").Value = 0.4239096552953754
$ws.Range("H23").Value = 0.4168677050772089
$ws.Range("H24").Value = 0.3905864257188227
$ws.Range("H25").Value = 0.3630617902832967
$ws.Range("I2").Value = 0.3757315010515043
$ws.Range("I3").Value = 0.3656238581545637
$ws.Range("I4").Value = 0.3596127975951688
$ws.Range("I5").Value = 0.3572128864359954
$ws.Range("I6").Value = 0.3568173988915788
$ws.Range("I7").Value = 0.3595802297560198
$ws.Range("I8").Value = 0.3722062484302811
$ws.Range("I9").Value = 0.3984893586955494
$ws.Range("I10").Value = 0.4186982167615625
$ws.Range("I11").Value = 0.4280806917633839
$ws.Range("I12").Value = 0.4316602591565939
$ws.Range("I13").Value = 0.4308881605543888
$ws.Range("I14").Value = 0.4283746547068787
$ws.Range("I15").Value = 0.4268385096919047
$ws.Range("I16").Value = 0.4180888088535895
$ws.Range("I17").Value = 0.4127692216154415
$ws.Range("I18").Value = 0.4097273953840173
$ws.Range("I19").Value = 0.4087005711172083
$ws.Range("I20").Value = 0.4133336567293782
$ws.Range("I21").Value = 0.429112214418673
$ws.Range("I22").Value = 0.4395793436413413
$ws.Range("I23").Value = 0.4339788619107594
$ws.Range("I24").Value = 0.4130784242511822
$ws.Range("I25").Value = 0.3912191043998661
$ws.Range("N2").Value = 1.538965158849294
$ws.Range("N3").Value = 1.440536212657179
$ws.Range("N4").Value = 1.380350891855613
$ws.Range("N5").Value = 1.355891405433653
$ws.Range("N6").Value = 1.351834048475041
$ws.Range("N7").Value = 1.38002074861987
$ws.Range("N8").Value = 1.504976964632107
$ws.Range("N9").Value = 1.751860150726543
$ws.Range("N10").Value = 1.93419659562025
$ws.Range("N11").Value = 2.017317443662478
$ws.Range("N12").Value = 2.048815004509549
$ws.Range("N13").Value = 2.042030543333738
$ws.Range("N14").Value = 2.019908357050213
$ws.Range("N15").Value = 2.006360570313291
$ws.Range("N16").Value = 1.928767661017872
$ws.Range("N17").Value = 1.88120914429345
$ws.Range("N18").Value = 1.85387144156806
$ws.Range("N19").Value = 1.844618345498787
$ws.Range("N20").Value = 1.886270128093656
$ws.Range("N21").Value = 2.026405625114307
$ws.Range("N22").Value = 2.118115898853375
$ws.Range("N23").Value = 2.069158297356239
$ws.Range("N24").Value = 1.883982043575145
$ws.Range("N25").Value = 1.684892220676318
$ws.Range("O2").Value = 1.288635611322178
$ws.Range("O3").Value = 1.186000399322808
$ws.Range("O4").Value = 1.123707417123967
$ws.Range("O5").Value = 1.098503372127993
$ws.Range("O6").Value = 1.094329137056548
$ws.Range("O7").Value = 1.12336677637748
$ws.Range("O8").Value = 1.253095441307636
$ws.Range("O9").Value = 1.513334097579957
$ws.Range("O10").Value = 1.70823048577472
$ws.Range("O11").Value = 1.797728296575315
$ws.Range("O12").Value = 1.831741466706262
$ws.Range("O13").Value = 1.824410655378244
$ws.Range("O14").Value = 1.800524119720649
$ws.Range("O15").Value = 1.785908904299276
$ws.Range("O16").Value = 1.702398612130764
$ws.Range("O17").Value = 1.651383722660285
$ws.Range("O18").Value = 1.622120123892898
$ws.Range("O19").Value = 1.612225456820568
$ws.Range("O20").Value = 1.656806175146926
$ws.Range("O21").Value = 1.807536841490844
$ws.Range("O22").Value = 1.906761666784462
$ws.Range("O23").Value = 1.853737639956933
$ws.Range("O24").Value = 1.654354480355494
$ws.Range("O25").Value = 1.442293097308493
